$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Nº hogares"
$ws.Range("B1").Value = "Tipo de hogar, código"
$ws.Range("C1").Value = "Municipio codigo"
$ws.Range("D1").Value = "Tipo de hogar"
$ws.Range("E1").Value = "Municipio nombre"

$ws.Range("A2").Value = "iaest-measure:n-hogares"
$ws.Range("B2").Value = "null"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "iaest-measure:tipo-de-hogar"
$ws.Range("E2").Value = "sdmx-dimension:refArea"

$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "null"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "dim"

$ws.Range("A4").Value = "xsd:double"
$ws.Range("B4").Value = "null"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "xsd:string"
$ws.Range("E4").Value = "URI-Municipio"
